# Smoke-test verification update for the regslv RTL testbench:
# the reset-signal column becomes a "synchronous reset signal" column,
# and each field row now documents which sync reset(s) it actually uses
# instead of the generic "Global Reset" placeholder.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column header: "复位信号" (Reset signal) -> "同步复位信号" (Sync reset signal).
# Also pick up the Chinese-header formatting used by the neighboring D8/E8
# cells (the old English-style header format no longer fits).
$ws.Range("G8").Value = "同步复位信号"
$ws.Range("D8").Copy()
$ws.Range("G8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Reserved bits rows have no associated reset net.
$ws.Range("G9").Value  = "None"
$ws.Range("G13").Value = "None"

# FIELD_1 / FIELD_2 now reference the concrete sync reset nets that drive them.
$ws.Range("G10").Value = "srst_10, srst_11"
$ws.Range("G11").Value = "srst_20"

# FIELD_3's read-user access was simplified to a plain read ("R"), and it
# no longer carries its own dedicated reset net.
$ws.Range("D12").Value = "R"
$ws.Range("G12").Value = "None"
